$d = $word.ActiveDocument

# Locate the run of text that needs to be split: " FinOps and Partnerships"
# (part of "SiriusXM | Director, FinOps and Partnerships") and change the
# job title to "SiriusXM | Director, Cloud FinOps".
$find = $d.Content.Find
$found = $find.Execute(" FinOps and Partnerships", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target text ' FinOps and Partnerships'"
}

$hit = $find.Parent
$startPos = $hit.Start
$endPos = $hit.End

# Re-fetch a fresh Range bound to the document (rather than reusing the
# Find's Parent range) so that InsertXML replaces this range's contents
# instead of merely appending after it.
$targetRange = $d.Range($startPos, $endPos)

# Replacement OOXML: keep a lone space (formatting identical to the
# original run), then two new bold runs "Cloud " and "FinOps" with the
# same run formatting as the text being replaced.
$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Gulim" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Gulim" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Cloud </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsia="Gulim" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>FinOps</w:t></w:r>' + `
'</w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)
